# The commit swaps the contents of ppt/theme/theme1.xml and ppt/theme/theme2.xml:
# theme1.xml (was "Office Theme") gets the "Integral" colour scheme and
# theme2.xml (was "Integral") gets the "Office Theme" colour scheme.
#
# ppt/theme/theme2.xml is the theme actually wired to the deck's single
# slide master (and to the presentation's primary theme relationship), so
# it is reachable from the PowerPoint object model via
# SlideMaster.Theme.ThemeColorScheme. We rewrite its 12 scheme colours to
# the values the "Office Theme" colour scheme (theme1.xml's current
# content) uses.

function ColorBgr($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# Order matches ThemeColorScheme.Colors(1..12): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = ColorBgr $officeThemeColors[$i - 1]
}
